# Update the cryptos list (Price / Volume(1h) columns) with the latest
# scraped values, per the "Updated cryptos list ... with GitHub Actions"
# commit. Only columns D (Price) and E (Volume(1h)) change; only for the
# rows whose figures actually moved.
#
# Note: columns D/E are stored as plain text (e.g. "27.872.75",
# "  -0.30%  ") rather than numbers, since values like "X.XXX.XX" aren't
# valid numerics and the "%" column keeps its literal padding/spacing.
# Assigning a plain numeric-looking string (e.g. "0.999") to .Value would
# be auto-coerced to a real number by Excel and could silently drop
# formatting such as trailing zeros (e.g. "3.10" -> 3.1). A leading
# apostrophe forces those values to stay text, exactly like the source
# data, while values that aren't valid numerics (e.g. "27.872.75", which
# has two dots) are assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.872.75'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.630.33'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '''211.62'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = '''0.518'
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '''0.0881'
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '1.860.61'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '1.640.09'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '''4.02'
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = '''0.557'
$ws.Range("E15").Value = '  -1.01%  '
$ws.Range("D16").Value = '''64.99'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("D17").Value = '27.893.97'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '''228.76'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").Value = '''0.998'
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").Value = '''4.35'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  -4.55%  '
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").Value = '''155.37'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = '''6.93'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '''15.47'
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").Value = '''0.0481'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("D33").Value = '1.418.56'
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").Value = '''3.10'
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("E35").Value = '  +2.32%  '
$ws.Range("D36").Value = '''1.01'
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("D37").Value = '''2.33'
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").Value = '''0.855'
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("D42").Value = '''65.93'
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").Value = '1.770.65'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("E46").Value = '  -3.82%  '
$ws.Range("D47").Value = '''88.71'
$ws.Range("E47").Value = '  +0.63%  '
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").Value = '''7.60'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Value = '''0.998'
$ws.Range("E51").Value = '  -0.26%  '
